$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for account 004346716 / TIAGO / 1604 first (it is below the
# other row to remove, so deleting it first keeps the other row index valid).
$ws.Rows.Item(23).Delete()

# Delete the row for account 005142624 / RODRIGO / 5000
$ws.Rows.Item(13).Delete()
